$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Location Name" column (E2:E9) with the new location values.
$ws.Range("E2").Value = "Bellagio Hotel"
$ws.Range("E3").Value = "The Getty"
$ws.Range("E4").Value = "Flatiron"
$ws.Range("E5").Value = "KOIN Center"
$ws.Range("E6").Value = "The Parthenon"
$ws.Range("E7").Value = "Olympia Theater"
$ws.Range("E8").Value = "Space Needle"
$ws.Range("E9").Value = "Coors Field"

# Move the active selection to E10, matching the saved cursor position.
$ws.Range("E10").Select()
